$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Swords & Daggers")
Write-Host $ws.Name
